$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Excuses")

# Remove the second logged excuse (row 3) entirely - shifts nothing below it up
$ws.Rows.Item(3).Delete()

# Student ID (A2) changed from 191480 to 190807.
# A leading apostrophe keeps it text (matching the original inline-string
# cell type) instead of Excel auto-coercing the digits to a number; then
# re-apply B2's cell format onto A2 so the quote-prefix flag set by the
# apostrophe entry doesn't linger as a stray style on the cell.
$ws.Range("A2").Value = "'190807"
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Log Date (C2) changed from 08/09/2025 to 20/09/2025
$ws.Range("C2").Value = "20/09/2025"
